$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Berenjena (Femacal de La Calera).
# It belongs right after the header/first data block, at row 29 - insert a
# whole row there (this shifts every existing record down by one, which is
# exactly what the workbook diff shows: old row N's data now lives in row
# N+1, all the way down to the former last row 139 becoming row 140).
$ws.Rows(29).Insert()

# Populate the newly inserted row 29 with the new record's values.
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44453
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 100112001
$ws.Range("G29").Value = "Berenjena"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 9000
$ws.Range("L29").Value = 9500
$ws.Range("M29").Value = 9200
$ws.Range("N29").Value = "$/caja 60 unidades"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 153
$ws.Range("Q29").Value = 60
$ws.Range("R29").Value = "Hortaliza"
